$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D range to text format while writing so Excel's smart
# number-detection does not reinterpret values like "663.025 €" as a
# float (663.025) with an implied currency unit; then restore the
# original (default/"Normal") cell style so formatting matches the source.
$dataRange = $ws.Range("D2:D26")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "1.599.048 €"
$ws.Range("D3").Value = "663.025 €"
$ws.Range("D4").Value = "502.512 €"
$ws.Range("D5").Value = "3.254.289 €"
$ws.Range("D6").Value = "3.911.408 €"
$ws.Range("D7").Value = "1.916.469 €"
$ws.Range("D8").Value = "768.794 €"
$ws.Range("D9").Value = "512.346 €"
$ws.Range("D10").Value = "3.288.733 €"
$ws.Range("D11").Value = "4.244.165 €"
$ws.Range("D12").Value = "2.617.352 €"
$ws.Range("D13").Value = "1.006.839 €"
$ws.Range("D14").Value = "784.491 €"
$ws.Range("D15").Value = "4.482.867 €"
$ws.Range("D16").Value = "5.618.742 €"
$ws.Range("D17").Value = "2.151.550 €"
$ws.Range("D18").Value = "792.234 €"
$ws.Range("D19").Value = "576.805 €"
$ws.Range("D20").Value = "4.090.277 €"
$ws.Range("D21").Value = "4.630.987 €"
$ws.Range("D22").Value = "2.667.921 €"
$ws.Range("D23").Value = "1.169.038 €"
$ws.Range("D24").Value = "617.064 €"
$ws.Range("D25").Value = "4.674.327 €"
$ws.Range("D26").Value = "5.889.525 €"

$dataRange.Style = "Normal"
